$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = @"
questions = [
    {
        "title": "Which commercial model (compensation structure) most requires a well-defined scope of work for the project team?",
        "ques_type": 2,
        "options": [
            "Hourly charge-out rates",
            "Fixed price",
            "Upfront retainer",
            "Contingent fee"
        ],
        "score": "Fixed price"
    },
    {
        "title": "You and your client are considering a contingent fee for your next project together, meaning that your compensation will depend on how successful the project is. What are the three key questions that need to be addressed to ensure a mutually-beneficial outcome?",
        "ques_type": 15,
        "options": [
            "How will the project\u2019s success be objectively measured?",
            "Will the potential fee closely resemble the \u201chourly rate\u201d fee you could earn otherwise?",
            "Will a partial fee be paid for a partially-successful project?",
            "How will success be defined?",
            "Will the fee be paid in a lump sum or in payments over time?",
            "Will there be a separate financial incentive for the project manager alone?",
            "Is this project\u2019s risk/reward ratio acceptable for both organizations?"
        ],
        "score": [
            "How will the project\u2019s success be objectively measured?",
            "How will success be defined?",
            "Is this project\u2019s risk/reward ratio acceptable for both organizations?"
        ]
    },
    {
        "title": "Your project is 50% complete, and there has been a three-week delay to a task that was not on the critical path and had two weeks of slack.  By how much should you expect the overall project to be delayed?",
        "ques_type": 2,
        "options": [
            "Zero weeks",
            "One week",
            "Three weeks",
            "Five weeks"
        ],
        "score": "One week"
    },
    {
        "title": "Today, you are presenting your results to your external clients and stakeholders from different departments within their firm whom you have not yet met. You submitted the presentation and the accompanying full report to your client last week they distributed both within their firm so that everyone could come prepared today with their questions. One of the less-influential stakeholders is monopolizing the question and answer period with basic questions, and you can sense your client\u2019s growing frustration. Your client is chairing the meeting. How should you handle the situation?",
        "ques_type": 2,
        "options": [
            "Point out the report section where this individual can find their answers.",
            "Suggest that this individual give others a turn to ask questions.",
            "Politely answer every question, no matter how basic.",
            "Reframe this individual\u2019s questions toward the key outcomes to benefit the group."
        ],
        "score": "Politely answer every question, no matter how basic."
    }
]
"@

# Here-strings include a trailing newline before the closing "@ marker; strip it off
$newText = $newText.TrimEnd("`r", "`n")

# The old layout had data in A1 (bold/bordered "0") and A2 (the shared string).
# The new layout collapses everything into A1 (plain/default style) and removes A2.
$ws.Range("A2").ClearContents()

# A1 previously used a bold font + thin border + centered/top alignment style;
# the new version uses the default (unstyled) cell format, so strip all formatting first.
$ws.Range("A1").ClearFormats()

$ws.Range("A1").Value = $newText

# Re-fit the row height to the new (multi-line) content instead of leaving a
# stale/explicit custom height behind.
$ws.Rows.Item(1).AutoFit()
